$d = $word.ActiveDocument

# Locate the last existing log entry ("Zelfde kleurcombinatie ...") so the
# insertion point is resilient to exact paragraph-index assumptions.
$anchorText = "Zelfde kleurcombinatie gebruiken voor login pagina"
$rng = $d.Content
$found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph '$anchorText'"
}

# Expand the found range to the whole paragraph (including its end-of-paragraph
# mark) and collapse to its end -- that's where the new log entries go, right
# before the trailing blank paragraph / sectPr.
[void]$rng.Expand(4)
[void]$rng.Collapse(0)

# Insert the new paragraphs (list items + trailing blank "Lijstalinea" para) as
# raw WordprocessingML via the standard single-part WordOpenXML package format.
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Index en login zijn gedeeld op </w:t></w:r><w:r><w:t>GitHub</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Begonnen aan magazijn pagina</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Product categorieën zijn toegevoegd, met voorbeeld categorieën</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Styling toegevoegd</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Lijsten verder uitwerken</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Alle lijsten zijn zichtbaar, en hebben hun styling gekregen</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>(Note voor morgen)</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Beginnen aan sorteerbuttons.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$rng.InsertXML($xml)

Write-Output "Inserted 8 new paragraphs after '$anchorText'."
